$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old "Total Duration:" summary that used to live in row 6 (cols C, D)
$ws.Range("C6").ClearContents()
$ws.Range("D6").ClearContents()

# Copy the existing data-row formatting onto the new cells before filling them in
$ws.Range("B7").Copy()
$ws.Range("C7:D8").PasteSpecial(-4122)

# Fill in the clock-out time and duration for the new entry in row 7
$ws.Range("C7").Value = "22:22:30"
$ws.Range("D7").Value = "-0.15 Hours"

# Add the new totals row at row 8
$ws.Range("C8").Value = "Total Duration:"
$ws.Range("D8").Value = "10.5 Hours"
